$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4800
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496
# Row 67
$ws.Range("H67").Value = 4800
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716
# Row 113
$ws.Range("H113").Value = 4090
$ws.Range("I113").Value = 3975
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 3975
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = -721
$ws.Range("N113").Value = -10674.6665

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1037.3334
$ws.Range("I2").Value = 1037.3334
$ws.Range("K2").Value = 1037.3334
$ws.Range("M2").Value = -924.3334
# Row 32
$ws.Range("H32").Value = 3502798.8
$ws.Range("I32").Value = 3336265.2
$ws.Range("K32").Value = 3336265.2
$ws.Range("M32").Value = -3335978.2
# Row 116
$ws.Range("H116").Value = 1037.3334
$ws.Range("I116").Value = 1037.3334
$ws.Range("K116").Value = 1037.3334
$ws.Range("M116").Value = 1256.6666

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1037.3334
$ws.Range("I3").Value = 1037.3334
$ws.Range("K3").Value = 1037.3334
$ws.Range("M3").Value = -923.3334
# Row 54
$ws.Range("H54").Value = 2796.6667
$ws.Range("J54").Value = 3500
$ws.Range("L54").Value = 3500
$ws.Range("N54").Value = -4468
# Row 86
$ws.Range("H86").Value = 1733.375
$ws.Range("I86").Value = 1733.375
$ws.Range("K86").Value = 1733.375
$ws.Range("M86").Value = -610.375
# Row 89
$ws.Range("H89").Value = 1733.375
$ws.Range("I89").Value = 1733.375
$ws.Range("K89").Value = 8666.875
$ws.Range("M89").Value = -3050.875
# Row 105
$ws.Range("H105").Value = 2231.6667
$ws.Range("I105").Value = 2231.6667
$ws.Range("K105").Value = 2231.6667
$ws.Range("M105").Value = -484.6667000000002
# Row 107
$ws.Range("H107").Value = 499.16666
$ws.Range("I107").Value = 498.33334
$ws.Range("K107").Value = 498.33334
$ws.Range("M107").Value = 1421.66666

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 5005
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5005
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5005
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -5283
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = $null
$ws.Range("N14").Value = $null
# Row 31
$ws.Range("H31").Value = 1654
$ws.Range("I31").Value = 1455.7142
$ws.Range("K31").Value = 1455.7142
$ws.Range("M31").Value = -1160.7142
# Row 34
$ws.Range("H34").Value = 1654
$ws.Range("I34").Value = 1455.7142
$ws.Range("K34").Value = 1455.7142
$ws.Range("M34").Value = -1253.7142
# Row 58
$ws.Range("H58").Value = 1819.8889
$ws.Range("I58").Value = 1292.2858
$ws.Range("K58").Value = 1292.2858
$ws.Range("M58").Value = -1089.2858
# Row 62
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
# Row 65
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
# Row 136
$ws.Range("H136").Value = 1819.8889
$ws.Range("I136").Value = 1292.2858
$ws.Range("K136").Value = 3876.8574
$ws.Range("M136").Value = -1326.8574
# Row 141
$ws.Range("H141").Value = 36824.8
$ws.Range("J141").Value = 36824.8
$ws.Range("L141").Value = 36824.8
$ws.Range("N141").Value = -47184.8

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 199
$ws.Range("J12").Value = 176.25
$ws.Range("L12").Value = 528.75
$ws.Range("N12").Value = -874.75
# Row 37
$ws.Range("H37").Value = 69946.336
$ws.Range("J37").Value = 69946.336
$ws.Range("L37").Value = 209839.008
$ws.Range("N37").Value = -210063.008
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = $null
# Row 76
$ws.Range("H76").Value = 16573.285
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 60000
$ws.Range("N76").Value = -60766
# Row 79
$ws.Range("H79").Value = 16573.285
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 60000
$ws.Range("N79").Value = -62652
# Row 111
$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567
# Row 120
$ws.Range("H120").Value = 3000
$ws.Range("I120").Value = 3000
$ws.Range("K120").Value = 9000
$ws.Range("M120").Value = -4162
# Row 139
$ws.Range("H139").Value = 2561
$ws.Range("I139").Value = 1149
$ws.Range("J139").Value = 3973
$ws.Range("K139").Value = 3447
$ws.Range("L139").Value = 11919
$ws.Range("M139").Value = 1693
$ws.Range("N139").Value = -22199

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 105460.8
$ws.Range("I10").Value = 172666.67
$ws.Range("J10").Value = 4652
$ws.Range("K10").Value = 172666.67
$ws.Range("L10").Value = 4652
$ws.Range("M10").Value = -172497.67
$ws.Range("N10").Value = -4990
# Row 102
$ws.Range("H102").Value = 1921.4546
$ws.Range("I102").Value = 1956.8422
$ws.Range("J102").Value = 1697.3334
$ws.Range("K102").Value = 1956.8422
$ws.Range("L102").Value = 1697.3334
$ws.Range("M102").Value = -334.8422
$ws.Range("N102").Value = -4941.3334
# Row 126
$ws.Range("H126").Value = 1928
$ws.Range("I126").Value = 1928
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5784
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3314
$ws.Range("N126").Value = $null

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 2552499.8
$ws.Range("I2").Value = 100000
$ws.Range("J2").Value = 5004999.5
$ws.Range("K2").Value = 100000
$ws.Range("L2").Value = 5004999.5
$ws.Range("M2").Value = -99888
$ws.Range("N2").Value = -5005223.5
# Row 16
$ws.Range("H16").Value = 359
$ws.Range("I16").Value = 359
$ws.Range("K16").Value = 359
$ws.Range("M16").Value = -189
# Row 22
$ws.Range("H22").Value = 596
$ws.Range("I22").Value = 719.6667
$ws.Range("K22").Value = 719.6667
$ws.Range("M22").Value = -424.6667
# Row 27
$ws.Range("H27").Value = 596
$ws.Range("I27").Value = 719.6667
$ws.Range("K27").Value = 719.6667
$ws.Range("M27").Value = -612.6667
# Row 40
$ws.Range("H40").Value = 3982
$ws.Range("I40").Value = 3646.8572
$ws.Range("K40").Value = 3646.8572
$ws.Range("M40").Value = -3510.8572
# Row 55
$ws.Range("H55").Value = 1365.8334
$ws.Range("I55").Value = 432.33334
$ws.Range("J55").Value = 2299.3333
$ws.Range("K55").Value = 432.33334
$ws.Range("L55").Value = 2299.3333
$ws.Range("M55").Value = -259.33334
$ws.Range("N55").Value = -2645.3333
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null
# Row 100
$ws.Range("H100").Value = 3838
$ws.Range("I100").Value = 3802.111
$ws.Range("J100").Value = 3999.5
$ws.Range("K100").Value = 3802.111
$ws.Range("L100").Value = 3999.5
$ws.Range("M100").Value = -3261.111
$ws.Range("N100").Value = -5081.5
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 375378.8
$ws.Range("J2").Value = 375378.8
$ws.Range("L2").Value = 375378.8
$ws.Range("N2").Value = -375602.8
# Row 41
$ws.Range("H41").Value = 19647.285
$ws.Range("I41").Value = 19668
$ws.Range("K41").Value = 19668
$ws.Range("M41").Value = -19278
# Row 54
$ws.Range("H54").Value = 44999.75
$ws.Range("J54").Value = 44999.75
$ws.Range("L54").Value = 44999.75
$ws.Range("N54").Value = -46039.75
# Row 70
$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50629
# Row 73
$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -52183
# Row 75
$ws.Range("H75").Value = 73003
$ws.Range("J75").Value = 73003
$ws.Range("L75").Value = 73003
$ws.Range("N75").Value = -74875
# Row 78
$ws.Range("H78").Value = 73003
$ws.Range("J78").Value = 73003
$ws.Range("L78").Value = 219009
$ws.Range("N78").Value = -228369
# Row 126
$ws.Range("H126").Value = 4399.3
$ws.Range("I126").Value = 4274.25
$ws.Range("K126").Value = 12822.75
$ws.Range("M126").Value = -10352.75
